# Add a new worksheet "Sheet2" after the existing "Sheet1", populate it with
# login-data rows (mirroring Sheet1's username/password layout but with a new
# "pass" column header and extra admin credentials), turn the email/password
# cells in row 2 into hyperlinks, and leave Sheet2 as the active sheet/tab.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# --- Cell data -------------------------------------------------------------
$ws2.Range("A1").Value = "username"
$ws2.Range("B1").Value = "pass"

$ws2.Range("A2").Value = "arun.joseph@learnship.com"
$ws2.Range("B2").Value = "Airtel@123"

$ws2.Range("A3").Value = "Insightadmin"
$ws2.Range("B3").Value = "Insight_0217"

# --- Column widths (approximate best-fit sizing) ----------------------------
$ws2.Columns("A").ColumnWidth = 34.25
$ws2.Columns("B").ColumnWidth = 10.5

# --- Hyperlinks on the credential row ---------------------------------------
$ws2.Hyperlinks.Add($ws2.Range("A2"), "mailto:arun.joseph@learnship.com")
$ws2.Hyperlinks.Add($ws2.Range("B2"), "mailto:Airtel@123")

# --- Selection / active sheet ------------------------------------------------
[void]$ws2.Range("D4").Select()
$ws2.Activate()
